# Replace three-digit division problems with new values.
$d = $word.ActiveDocument

$replacements = @(
    @("950÷7=135, 5", "389÷3=129, 2"),
    @("154÷3=51, 1", "779÷9=86, 5"),
    @("721÷9=80, 1", "788÷6=131, 2"),
    @("721÷8=90, 1", "843÷3=281, 0"),
    @("635÷5=127, 0", "397÷7=56, 5"),
    @("270÷5=54, 0", "926÷4=231, 2"),
    @("227÷6=37, 5", "946÷3=315, 1"),
    @("388÷2=194, 0", "820÷4=205, 0"),
    @("925÷5=185, 0", "301÷4=75, 1"),
    @("326÷4=81, 2", "127÷4=31, 3"),
    @("180÷2=90, 0", "989÷7=141, 2"),
    @("169÷2=84, 1", "922÷7=131, 5"),
    @("404÷6=67, 2", "269÷3=89, 2"),
    @("940÷3=313, 1", "648÷7=92, 4"),
    @("397÷5=79, 2", "542÷2=271, 0"),
    @("480÷8=60, 0", "577÷3=192, 1"),
    @("780÷4=195, 0", "261÷9=29, 0"),
    @("487÷7=69, 4", "962÷2=481, 0"),
    @("131÷5=26, 1", "783÷2=391, 1"),
    @("226÷3=75, 1", "895÷7=127, 6"),
    @("929÷4=232, 1", "997÷8=124, 5"),
    @("954÷9=106, 0", "724÷3=241, 1"),
    @("553÷3=184, 1", "476÷8=59, 4"),
    @("872÷5=174, 2", "517÷7=73, 6"),
    @("126÷7=18, 0", "245÷5=49, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
